$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 128 (shifts existing rows 128:137 down to 129:138)
$ws.Rows.Item(128).Insert()

# Populate the new row 128 with the new weekly data point
$ws.Range("A128").Value = 9
$ws.Range("B128").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C128").Value = "Metropolitana"
$ws.Range("D128").Value = 45131
$ws.Range("E128").Value = 13
$ws.Range("F128").Value = 100114007
$ws.Range("G128").Value = "Jengibre"
$ws.Range("H128").Value = "Sin especificar"
$ws.Range("I128").Value = "Primera"
$ws.Range("J128").Value = 520
$ws.Range("K128").Value = 17000
$ws.Range("L128").Value = 18000
$ws.Range("M128").Value = 17500
$ws.Range("N128").Value = "`$/caja 13 kilos"
$ws.Range("O128").Value = "Perú"
$ws.Range("P128").Value = 1346
$ws.Range("Q128").Value = 13
$ws.Range("R128").Value = "Hortaliza"
